# Generate Report for Handback
#
# Overview sheet: both zh-cn and de-de status cells move from
# "Ready for handoff" to "Handed back: in sync with en-US".
#
# zh-cn / de-de sheets: the handback round-trip is recorded -
# "Latest Target File" (I2) and "Latest Handback File" (J2) get filled in,
# and "Latest Handback DateTime" (K2) is stamped with the handback time.

$wb = $excel.ActiveWorkbook

$targetFileName   = "a04b0430-406c-4e47-9e35-46b35a874fe4.md"
$targetFileUrl    = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cc43797710dabeb9e82def370c1eb0652317f6d4/e2e/a04b0430-406c-4e47-9e35-46b35a874fe4.md"
$zhHandbackXlf    = "a04b0430-406c-4e47-9e35-46b35a874fe4.3914c112f2ce338a9c100e4757515d7effef857b.zh-cn.xlf"
$deHandbackXlf    = "a04b0430-406c-4e47-9e35-46b35a874fe4.3914c112f2ce338a9c100e4757515d7effef857b.de-de.xlf"
$newStatus        = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet - update the per-language status cells
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Columns.Item(5).AutoFit()
$overview.Columns.Item(6).AutoFit()

# ---------------------------------------------------------------------
# zh-cn sheet - status + handback columns
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Columns.Item(3).AutoFit()

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $targetFileUrl, "", "", $targetFileName)
$zhcn.Range("I2").Style = "HyperLink"
$zhcn.Range("J2").Value = $zhHandbackXlf
$zhcn.Range("K2").Value = "2016-08-17 08:57:57"
$zhcn.Columns.Item(9).AutoFit()
$zhcn.Columns.Item(10).AutoFit()

# ---------------------------------------------------------------------
# de-de sheet - status + handback columns
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Columns.Item(3).AutoFit()

$dede.Hyperlinks.Add($dede.Range("I2"), $targetFileUrl, "", "", $targetFileName)
$dede.Range("I2").Style = "HyperLink"
$dede.Range("J2").Value = $deHandbackXlf
$dede.Range("K2").Value = "2016-08-17 08:58:13"
$dede.Columns.Item(9).AutoFit()
$dede.Columns.Item(10).AutoFit()
